$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $A, $B, $C, $D, $E, $F, $G)
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
}

Set-RowData 86 2026 "ENERO" 4 "AMARILLO" "COLORES" "GOLDFINCH" 10052
Set-RowData 87 2026 "ENERO" 4 "AMARILLO" "COLORES" "HIGH AND EXOTIC" 26771
Set-RowData 88 2026 "ENERO" 4 "AMARILLO" "COLORES" "MOMENTUM" 12450
Set-RowData 89 2026 "ENERO" 4 "AMARILLO" "COLORES" "SUNDAY MORNING" 1003
Set-RowData 90 2026 "ENERO" 4 "BIC. AMARILLO" "COLORES" "SUMMER LIGHT" 6518
Set-RowData 91 2026 "ENERO" 4 "BICOLOR" "COLORES" "BLUSH" 2928
Set-RowData 92 2026 "ENERO" 4 "BICOLOR" "COLORES" "DISCOVERY" 850
Set-RowData 93 2026 "ENERO" 4 "BLANCO" "COLORES" "HIGH AND PURE" 24185
Set-RowData 94 2026 "ENERO" 4 "BLANCO" "COLORES" "SUGGAR DOLL" 7860
Set-RowData 95 2026 "ENERO" 4 "BLANCO" "COLORES" "VANILLA ICE" 7100
Set-RowData 96 2026 "ENERO" 4 "BLANCO" "COLORES" "VENDELA" 9950
Set-RowData 97 2026 "ENERO" 4 "DURAZNO" "COLORES" "TIFANY" 4560
Set-RowData 98 2026 "ENERO" 4 "HOT PINK" "COLORES" "COTTON CANDY" 7400
Set-RowData 99 2026 "ENERO" 4 "HOT PINK" "COLORES" "JACARANDA" 25700
Set-RowData 100 2026 "ENERO" 4 "HOT PINK" "COLORES" "PINK FLOYD" 13953
Set-RowData 101 2026 "ENERO" 4 "LAVANDER" "COLORES" "DEEP PURPLE" 32605
Set-RowData 102 2026 "ENERO" 4 "LAVANDER" "COLORES" "MOODY BLUES" 17226
Set-RowData 103 2026 "ENERO" 4 "NARANJA" "COLORES" "ALIVE" 6968
Set-RowData 104 2026 "ENERO" 4 "NARANJA" "COLORES" "BROMO" 8006
Set-RowData 105 2026 "ENERO" 4 "NARANJA" "COLORES" "CLEMENTINA" 3720
Set-RowData 106 2026 "ENERO" 4 "NARANJA" "COLORES" "NINA" 29103
Set-RowData 107 2026 "ENERO" 4 "ROJO" "ROJO" "FREEDOM" 304105
Set-RowData 108 2026 "ENERO" 4 "ROSADO" "COLORES" "ABSOLUT IN PINK" 7715
Set-RowData 109 2026 "ENERO" 4 "ROSADO" "COLORES" "HIGH AND BONITA" 23585
Set-RowData 110 2026 "ENERO" 4 "ROSADO" "COLORES" "LUCIANO" 2150
Set-RowData 111 2026 "ENERO" 4 "ROSADO" "COLORES" "PECKOUBO" 2550
Set-RowData 112 2026 "ENERO" 4 "ROSADO" "COLORES" "STARFISH" 4990
Set-RowData 113 2026 "ENERO" 4 "ROSADO" "COLORES" "TABATHA" 15989

[void]$ws.Range("A113").Select()
